# Replace Swaraj Tractors model/image data with Sonalika Tractors data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$brand = "Sonalika Tractors"

$models = @(
    "DI 35",
    "745 DI III Sikander",
    "745 DI III Sikander",
    "Tiger 50",
    "DI 750III",
    "Tiger Electric 4WD",
    "DI 60",
    "DI 42 RX",
    "DI 745 III",
    "42 DI Sikander",
    "WT 60"
)

$images = @(
    "['DI35img0-di-35-1631526570.png', 'DI35img1-mqdefault.png', 'DI35img2-mqdefault.png']",
    "['745DIIIISikanderimg0-745-di-iii-sikander-1631526194.png', '745DIIIISikanderimg1-upload-1631526194-0.png', '745DIIIISikanderimg2-upload-1631526194-0.png']",
    "['745DIIIISikanderimg0-745-di-iii-sikander-1631526194.png', '745DIIIISikanderimg1-upload-1631526194-0.png', '745DIIIISikanderimg2-upload-1631526194-0.png']",
    "['Tiger50img0-tiger-50-1631530291.png', 'Tiger50img1-upload-1631530291-0.png', 'Tiger50img2-upload-1631530291-0.png']",
    "['DI750IIIimg0-di-750iii-1631529270.png', 'DI750IIIimg1-upload-1631529270-0.png', 'DI750IIIimg2-di-750iii-1631529270.png']",
    "['TigerElectric4WDimg0-sonalika-tiger-electric-4wd-1696505405.png', 'TigerElectric4WDimg1-sonalika-tiger-electric-4wd-16965054050.png', 'TigerElectric4WDimg2-sonalika-tiger-electric-4wd-1696505405.png']",
    "['DI60img0-di-60-1631527166.png', 'DI60img1-upload-1631527166-0.png', 'DI60img2-upload-1631527166-0.png']",
    "['DI42RXimg0-di-42-rx-1631526666.png', 'DI42RXimg1-di-42-rx-1631526666.png', 'DI42RXimg2-mqdefault.png']",
    "['DI745IIIimg0-di-745-iii-1631528823.png', 'DI745IIIimg1-upload-1631528823-0.png', 'DI745IIIimg2-mqdefault.png']",
    "['42DISikanderimg0-42-di-sikander-1631525877.png', '42DISikanderimg1-upload-1631525877-0.png', '42DISikanderimg2-upload-1631525877-0.png']",
    "['WT60img0-wt-60-1631530732.png', 'WT60img1-upload-1631530732-0.png', 'WT60img2-wt-60-1631530732.png']"
)

for ($i = 0; $i -lt $models.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $brand
    $ws.Cells.Item($row, 2).Value = $models[$i]
    $ws.Cells.Item($row, 3).Value = $images[$i]
}
